# Update countries & provincias Spain
# Applies the data refresh captured in the commit: updated "last updated" timestamp,
# refreshed case counters for several countries, and two countries that swapped
# ranking order (so their row positions/labels swap together with their figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Septiembre de 2020 a las 07:54"

# --- Row 5: India ---
$ws.Range("B5").Value = 5562663
$ws.Range("C5").Value = 2558
$ws.Range("D5").Value = 4497867
$ws.Range("E5").Value = 975831

# --- Row 28: Ucrania ---
$ws.Range("B28").Value = 181237
$ws.Range("C28").Value = 2884
$ws.Range("D28").Value = 79901
$ws.Range("E28").Value = 97694
$ws.Range("G28").Value = 59
$ws.Range("H28").Value = 3642

# --- Row 59: Uzbekistan ---
$ws.Range("B59").Value = 52491
$ws.Range("C59").Value = 421
$ws.Range("D59").Value = 48606
$ws.Range("E59").Value = 3445
$ws.Range("G59").Value = 3
$ws.Range("H59").Value = 440

# --- Row 66: Kirguistan ---
$ws.Range("B66").Value = 45542
$ws.Range("C66").Value = 71
$ws.Range("D66").Value = 41796
$ws.Range("E66").Value = 2683

# --- Row 77: El Salvador ---
$ws.Range("D77").Value = 21795
$ws.Range("E77").Value = 5189
$ws.Range("G77").Value = 2
$ws.Range("H77").Value = 814

# --- Rows 163/164: Liberia and Polinesia Francesa swap order ---
$ws.Range("A163").Value = "Polinesia Francesa"
$ws.Range("B163").Value = 1394
$ws.Range("C163").Value = 123
$ws.Range("D163").Value = 1170
$ws.Range("E163").Value = 222
$ws.Range("H163").Value = 2

$ws.Range("A164").Value = "Liberia"
$ws.Range("B164").Value = 1336
$ws.Range("D164").Value = 1218
$ws.Range("E164").Value = 36
$ws.Range("H164").Value = 82

# --- Rows 214/215: Montserrat and Islas Malvinas swap order ---
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1
